$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (style already bold/bordered from existing sheet formatting)
$ws.Range("A1").Value = "topic_id"
$ws.Range("B1").Value = "title"
$ws.Range("C1").Value = "user_id"
$ws.Range("D1").Value = "nickname"
$ws.Range("E1").Value = "content"

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "улгту"
$ws.Range("C2").Value = 5813154625
$ws.Range("D2").Value = "dfgdf"
$ws.Range("E2").Value = "привет"

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "лабы"
$ws.Range("C3").Value = 5813154625
$ws.Range("D3").Value = "dfgdf"
$ws.Range("E3").Value = "пока"

# Remove now-unused column F entirely (sheet shrinks from F to E)
$ws.Columns.Item(6).Delete()
